# Update "想去人数" (column F) counters across the four sheets to match
# the regenerated GitHub Pages data snapshot (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 447
$ws.Range("F5").Value = 1342
$ws.Range("F6").Value = 7684
$ws.Range("F7").Value = 96
$ws.Range("F9").Value = 2099
$ws.Range("F10").Value = 8472
$ws.Range("F13").Value = 67
$ws.Range("F14").Value = 5673
$ws.Range("F16").Value = 2631
$ws.Range("F17").Value = 1145
$ws.Range("F19").Value = 346
$ws.Range("F20").Value = 406
$ws.Range("F23").Value = 539
$ws.Range("F24").Value = 3576
$ws.Range("F26").Value = 39
$ws.Range("F27").Value = 28
$ws.Range("F29").Value = 3076
$ws.Range("F30").Value = 47
$ws.Range("F31").Value = 127
$ws.Range("F32").Value = 353
$ws.Range("F33").Value = 130
$ws.Range("F34").Value = 320
$ws.Range("F35").Value = 648
$ws.Range("F39").Value = 1950
$ws.Range("F43").Value = 2993
$ws.Range("F45").Value = 2294
$ws.Range("F49").Value = 2

# 演出 (Performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 130
$ws.Range("F9").Value = 124

# 本地生活 (Local life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 1333

# 全部类型 (All types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1333
$ws.Range("F5").Value = 1342
$ws.Range("F6").Value = 7684
$ws.Range("F7").Value = 96
$ws.Range("F9").Value = 2099
$ws.Range("F10").Value = 8472
$ws.Range("F12").Value = 67
$ws.Range("F13").Value = 5673
$ws.Range("F15").Value = 2631
$ws.Range("F16").Value = 1145
$ws.Range("F18").Value = 406
$ws.Range("F22").Value = 130
$ws.Range("F23").Value = 539
$ws.Range("F25").Value = 3576
$ws.Range("F27").Value = 39
$ws.Range("F28").Value = 28
$ws.Range("F30").Value = 3076
$ws.Range("F31").Value = 353
$ws.Range("F32").Value = 130
$ws.Range("F33").Value = 320
$ws.Range("F35").Value = 648
$ws.Range("F40").Value = 1951
$ws.Range("F44").Value = 2994
$ws.Range("F45").Value = 2294
$ws.Range("F48").Value = 124
